$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "specializations" column header in J1
$ws.Range("J1").Value = "specializations"

# Update the active selection to match the edited state
$ws.Range("I6").Select()
